$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the second column of data (currently in column D) to column F,
# leaving columns B, C, and E empty, per the new CSV layout for the 20q data.
for ($r = 1; $r -le 9; $r++) {
    $srcCell = $ws.Cells.Item($r, 4)   # column D
    $dstCell = $ws.Cells.Item($r, 6)   # column F
    $dstCell.Value = $srcCell.Value2
    $srcCell.ClearContents()
}

# Update the selected / active cell to D4
$ws.Range("D4").Select()
